$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.714.94"
$ws.Range("E2").Value = "'  +2.27%  "
$ws.Range("D3").Value = "'2.164.60"
$ws.Range("E3").Value = "'  +2.85%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'227.85"
$ws.Range("E5").Value = "'  +0.08%  "
$ws.Range("D6").Value = "'0.630"
$ws.Range("E6").Value = "'  +2.25%  "
$ws.Range("D7").Value = "'63.73"
$ws.Range("E7").Value = "'  +2.27%  "
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E9").Value = "'  +0.86%  "
$ws.Range("D10").Value = "'0.0848"
$ws.Range("E10").Value = "'  +0.55%  "
$ws.Range("E11").Value = "'  +0.14%  "
$ws.Range("D12").Value = "'16.01"
$ws.Range("E12").Value = "'  +1.23%  "
$ws.Range("D13").Value = "'2.483.11"
$ws.Range("E13").Value = "'  +2.69%  "
$ws.Range("D14").Value = "'21.99"
$ws.Range("E14").Value = "'  -0.39%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("E15").Value = "'  +0.14%  "
$ws.Range("D16").Value = "'5.50"
$ws.Range("E16").Value = "'  -0.59%  "
$ws.Range("D17").Value = "'2.178.04"
$ws.Range("E17").Value = "'  +4.14%  "
$ws.Range("D18").Value = "'39.625.09"
$ws.Range("E18").Value = "'  +1.98%  "
$ws.Range("D19").Value = "'71.83"
$ws.Range("E19").Value = "'  +0.24%  "
$ws.Range("E20").Value = "'  -0.32%  "
$ws.Range("D21").Value = "'0.0₃0846"
$ws.Range("E21").Value = "'  -0.07%  "
$ws.Range("D22").Value = "'228.43"
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("E23").Value = "'  +0.09%  "
$ws.Range("E24").Value = "'  +3.27%  "
$ws.Range("E25").Value = "'  -0.02%  "
$ws.Range("D26").Value = "'172.91"
$ws.Range("E26").Value = "'  +0.48%  "
$ws.Range("D27").Value = "'9.65"
$ws.Range("E27").Value = "'  -0.15%  "
$ws.Range("D29").Value = "'19.82"
$ws.Range("E29").Value = "'  +2.44%  "
$ws.Range("D30").Value = "'1.42"
$ws.Range("E30").Value = "'  +0.41%  "
$ws.Range("E31").Value = "'  +4.56%  "
$ws.Range("E32").Value = "'  +1.54%  "
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("E34").Value = "'  -1.34%  "
$ws.Range("D35").Value = "'6.96"
$ws.Range("E35").Value = "'  -3.03%  "
$ws.Range("D36").Value = "'0.0618"
$ws.Range("E36").Value = "'  -0.08%  "
$ws.Range("D37").Value = "'2.41"
$ws.Range("E37").Value = "'  +0.58%  "
$ws.Range("D38").Value = "'3.61"
$ws.Range("E38").Value = "'  +2.39%  "
$ws.Range("D39").Value = "'5.20"
$ws.Range("E39").Value = "'  +24.82%  "
$ws.Range("E40").Value = "'  -0.18%  "
$ws.Range("D41").Value = "'102.33"
$ws.Range("E41").Value = "'  +0.19%  "
$ws.Range("D42").Value = "'0.0227"
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("E43").Value = "'  +3.16%  "
$ws.Range("D44").Value = "'17.46"
$ws.Range("E44").Value = "'  -3.63%  "
$ws.Range("D45").Value = "'1.516.40"
$ws.Range("E45").Value = "'  -0.59%  "
$ws.Range("E46").Value = "'  +0.83%  "
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("D48").Value = "'1.10"
$ws.Range("E48").Value = "'  +1.35%  "
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = "'  +0.10%  "
$ws.Range("E50").Value = "'  +1.02%  "
$ws.Range("D51").Value = "'2.368.67"
$ws.Range("E51").Value = "'  +2.80%  "
